$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a literal text value into a cell so it keeps the workbook's
# original text cell-type (matching the source inlineStr cells) instead of
# Excel auto-coercing numeric-looking strings (e.g. "6.25") into numbers.
# For values that look numeric we briefly force a Text number format, assign
# the value, then restore the cell to the default "Normal" style so no stray
# formatting is left behind.
function Set-TextValue {
    param($cell, [string]$text)
    $looksNumeric = $text -match '^[+-]?[0-9]+(\.[0-9]+)?$'
    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Cells.Item(2, 4) '56.666.70'
Set-TextValue $ws.Cells.Item(2, 5) '  +0.57%  '
Set-TextValue $ws.Cells.Item(3, 4) '2.387.32'
Set-TextValue $ws.Cells.Item(3, 5) '  +0.81%  '
Set-TextValue $ws.Cells.Item(4, 5) '  -0.02%  '
Set-TextValue $ws.Cells.Item(5, 4) '504.26'
Set-TextValue $ws.Cells.Item(5, 5) '  +1.21%  '
Set-TextValue $ws.Cells.Item(6, 4) '132.53'
Set-TextValue $ws.Cells.Item(6, 5) '  +2.98%  '
Set-TextValue $ws.Cells.Item(7, 5) '  +0.13%  '
Set-TextValue $ws.Cells.Item(8, 4) '0.551'
Set-TextValue $ws.Cells.Item(8, 5) '  +0.23%  '
Set-TextValue $ws.Cells.Item(9, 4) '2.392.09'
Set-TextValue $ws.Cells.Item(9, 5) '  +0.13%  '
Set-TextValue $ws.Cells.Item(10, 5) '  +2.03%  '
Set-TextValue $ws.Cells.Item(11, 5) '  +0.72%  '
Set-TextValue $ws.Cells.Item(12, 5) '  +2.12%  '
Set-TextValue $ws.Cells.Item(13, 5) '  +1.56%  '
Set-TextValue $ws.Cells.Item(14, 4) '2.811.18'
Set-TextValue $ws.Cells.Item(14, 5) '  +0.63%  '
Set-TextValue $ws.Cells.Item(15, 4) '56.599.40'
Set-TextValue $ws.Cells.Item(15, 5) '  -0.44%  '
Set-TextValue $ws.Cells.Item(16, 4) '21.64'
Set-TextValue $ws.Cells.Item(16, 5) '  +1.04%  '
Set-TextValue $ws.Cells.Item(17, 5) '  +1.72%  '
Set-TextValue $ws.Cells.Item(18, 4) '2.381.02'
Set-TextValue $ws.Cells.Item(18, 5) '  -0.97%  '
Set-TextValue $ws.Cells.Item(19, 5) '  +0.73%  '
Set-TextValue $ws.Cells.Item(20, 5) '  +0.79%  '
Set-TextValue $ws.Cells.Item(21, 4) '308.76'
Set-TextValue $ws.Cells.Item(21, 5) '  -0.81%  '
Set-TextValue $ws.Cells.Item(22, 4) '6.25'
Set-TextValue $ws.Cells.Item(22, 5) '  +1.44%  '
Set-TextValue $ws.Cells.Item(23, 5) '  +0.15%  '
Set-TextValue $ws.Cells.Item(24, 5) '  -4.15%  '
Set-TextValue $ws.Cells.Item(25, 4) '65.97'
Set-TextValue $ws.Cells.Item(25, 5) '  +0.92%  '
Set-TextValue $ws.Cells.Item(26, 4) '0.996'
Set-TextValue $ws.Cells.Item(26, 5) '  -0.73%  '
Set-TextValue $ws.Cells.Item(27, 4) '0.381'
Set-TextValue $ws.Cells.Item(27, 5) '  +3.07%  '
Set-TextValue $ws.Cells.Item(28, 5) '  +0.60%  '
Set-TextValue $ws.Cells.Item(29, 4) '7.35'
Set-TextValue $ws.Cells.Item(29, 5) '  +2.63%  '
Set-TextValue $ws.Cells.Item(30, 4) '176.12'
Set-TextValue $ws.Cells.Item(30, 5) '  +1.24%  '
Set-TextValue $ws.Cells.Item(31, 4) '0.0₃0725'
Set-TextValue $ws.Cells.Item(31, 5) '  +2.76%  '
Set-TextValue $ws.Cells.Item(32, 5) '  +0.25%  '
Set-TextValue $ws.Cells.Item(33, 5) '  +2.20%  '
Set-TextValue $ws.Cells.Item(34, 4) '5.85'
Set-TextValue $ws.Cells.Item(34, 5) '  -3.51%  '
Set-TextValue $ws.Cells.Item(35, 5) '  +0.10%  '
Set-TextValue $ws.Cells.Item(36, 5) '  +0.21%  '
Set-TextValue $ws.Cells.Item(37, 5) '  +0.42%  '
Set-TextValue $ws.Cells.Item(38, 5) '  -2.04%  '
Set-TextValue $ws.Cells.Item(39, 5) '  +2.15%  '
Set-TextValue $ws.Cells.Item(40, 4) '36.76'
Set-TextValue $ws.Cells.Item(40, 5) '  +2.60%  '
Set-TextValue $ws.Cells.Item(41, 5) '  +6.98%  '
Set-TextValue $ws.Cells.Item(42, 4) '1.43'
Set-TextValue $ws.Cells.Item(42, 5) '  +1.25%  '
Set-TextValue $ws.Cells.Item(43, 4) '130.73'
Set-TextValue $ws.Cells.Item(43, 5) '  +1.41%  '
Set-TextValue $ws.Cells.Item(44, 4) '3.37'
Set-TextValue $ws.Cells.Item(44, 5) '  +1.11%  '
Set-TextValue $ws.Cells.Item(45, 4) '4.83'
Set-TextValue $ws.Cells.Item(45, 5) '  +1.61%  '
Set-TextValue $ws.Cells.Item(46, 5) '  -0.23%  '
Set-TextValue $ws.Cells.Item(47, 5) '  +1.37%  '
Set-TextValue $ws.Cells.Item(48, 4) '247.48'
Set-TextValue $ws.Cells.Item(48, 5) '  -1.86%  '
Set-TextValue $ws.Cells.Item(49, 5) '  +0.07%  '
Set-TextValue $ws.Cells.Item(50, 5) '  +1.95%  '
Set-TextValue $ws.Cells.Item(51, 4) '17.17'
Set-TextValue $ws.Cells.Item(51, 5) '  +7.85%  '
